$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Reason of change text for row 2 (Amietia hymenopus)
$ws.Range("D2").Value = "Deterioration in habitat quality, especially in Lesotho where active and proposed mines  have had an observable increased impact (both directly and indirectly) on this species, with projected increased aridification due to climate change, overgrazing and water quality issues that have intensified since the last assessment, resulting in this species’ conservation status deteriorating."

# Update the Reason of change text for row 3 (Breviceps macrops)
$ws.Range("D3").Value = "Based on projected increased threats from mining and development. New information on future threatening processes has become available with changes in government policy that will likely lead to increased industrial development in the coastal region of the Northern Cape, South Africa and southwestern Namibia."

# Remove row 4 (Hyperolius poweri) entirely - species entry removed from the table
$ws.Rows("4:4").Delete()

# Update row height to match the shorter wrapped text in row 3
$ws.Rows("3:3").RowHeight = 158.4

# Match the final selection / scroll position left by the author
$wb.Windows.Item(1).ScrollRow = 2
$ws.Range("D3").Select()
